# Atualiza parcial e avisos da Liga Eliminacao 20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header in column V, mirroring the style of the existing header row (A1:U1)
$ws.Range("V1").Value = "Parcial Rodada 1"
$ws.Range("A1").Copy()
$ws.Range("V1").PasteSpecial(-4122)  # xlPasteFormats

# Partial-round values for column V (rows 2-21)
$parciais = @(44.26, 40.4, 45.46, 58.17, 59.86, 16.4, 38.46, 39.66, 55.96, 63.76, 83.5, 0, 67.16, 53.06, 44.65, 62.76, 51.5, 52.66, 54.36, 53.91)

for ($i = 0; $i -lt $parciais.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 22).Value = $parciais[$i]
}
